$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Groups"
$ws.Range("C2").Value = "A"
$ws.Range("C3").Value = "A"
$ws.Range("C4").Value = "A"
$ws.Range("C5").Value = "B"
$ws.Range("C6").Value = "B"

$ws.Range("C7").Select()
